# Auto-generated script applying cell-level numeric updates to the
# "Sephirot_Profits" workbook sheets, per the scheduled-runner commit.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6793.6
$ws.Range("I64").Value = 6490.5
$ws.Range("K64").Value = 6490.5
$ws.Range("M64").Value = -6242.5
$ws.Range("H67").Value = 6793.6
$ws.Range("I67").Value = 6490.5
$ws.Range("K67").Value = 6490.5
$ws.Range("M67").Value = -5632.5
$ws.Range("H99").Value = 699.6667
$ws.Range("I99").Value = 299.5
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 898.5
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = 599.5
$ws.Range("N99").Value = -7496
$ws.Range("H111").Value = 491.6
$ws.Range("I111").Value = 531.5
$ws.Range("K111").Value = 1594.5
$ws.Range("M111").Value = 1472.5
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 565.5
$ws.Range("I2").Value = 578.6
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 578.6
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -465.6
$ws.Range("N2").Value = -726
$ws.Range("H116").Value = 565.5
$ws.Range("I116").Value = 578.6
$ws.Range("J116").Value = 500
$ws.Range("K116").Value = 578.6
$ws.Range("L116").Value = 500
$ws.Range("M116").Value = 1715.4
$ws.Range("N116").Value = -5088
$ws.Range("H132").Value = 1910.2368
$ws.Range("I132").Value = 1118.2593
$ws.Range("J132").Value = 3854.182
$ws.Range("K132").Value = 3354.7779
$ws.Range("L132").Value = 11562.546
$ws.Range("M132").Value = -824.7779
$ws.Range("N132").Value = -16622.546

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 565.5
$ws.Range("I3").Value = 578.6
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 578.6
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = -464.6
$ws.Range("N3").Value = -728
$ws.Range("H94").Value = 1724.75
$ws.Range("I94").Value = 1724.75
$ws.Range("K94").Value = 1724.75
$ws.Range("M94").Value = -1273.75
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 622.1
$ws.Range("I107").Value = 622.1
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 622.1
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1297.9
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 6642.143
$ws.Range("I134").Value = 1077
$ws.Range("K134").Value = 3231
$ws.Range("M134").Value = -696

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3249.5
$ws.Range("I16").Value = 2999
$ws.Range("K16").Value = 2999
$ws.Range("M16").Value = -2712
$ws.Range("H31").Value = 2364.8333
$ws.Range("I31").Value = 1130.6666
$ws.Range("K31").Value = 1130.6666
$ws.Range("M31").Value = -835.6666
$ws.Range("H34").Value = 2364.8333
$ws.Range("I34").Value = 1130.6666
$ws.Range("K34").Value = 1130.6666
$ws.Range("M34").Value = -928.6666
$ws.Range("H58").Value = 1000
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 1000
$ws.Range("M58").Value = -797
$ws.Range("H86").Value = 1999.5
$ws.Range("J86").Value = 1999
$ws.Range("L86").Value = 1999
$ws.Range("N86").Value = -4245
$ws.Range("H89").Value = 1999.5
$ws.Range("J89").Value = 1999
$ws.Range("L89").Value = 9995
$ws.Range("N89").Value = -21227
$ws.Range("H94").Value = 1500
$ws.Range("J94").Value = 1500
$ws.Range("L94").Value = 1500
$ws.Range("N94").Value = -2402
$ws.Range("H113").Value = 3249.5
$ws.Range("I113").Value = 2999
$ws.Range("K113").Value = 2999
$ws.Range("M113").Value = -829
$ws.Range("H122").Value = 1601.3334
$ws.Range("I122").Value = 959.6
$ws.Range("K122").Value = 2878.8
$ws.Range("M122").Value = -428.8000000000002
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 200
$ws.Range("J68").Value = 200
$ws.Range("L68").Value = 600
$ws.Range("N68").Value = -2222
$ws.Range("H71").Value = 200
$ws.Range("J71").Value = 200
$ws.Range("L71").Value = 1800
$ws.Range("N71").Value = -9912
$ws.Range("H113").Value = 320.55554
$ws.Range("I113").Value = 683.3333
$ws.Range("J113").Value = 248
$ws.Range("K113").Value = 2049.9999
$ws.Range("L113").Value = 744
$ws.Range("M113").Value = 120.0001000000002
$ws.Range("N113").Value = -5084

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3175.5
$ws.Range("I80").Value = 3163.5
$ws.Range("K80").Value = 3163.5
$ws.Range("M80").Value = -2165.5
$ws.Range("H83").Value = 3175.5
$ws.Range("I83").Value = 3163.5
$ws.Range("K83").Value = 15817.5
$ws.Range("M83").Value = -10825.5
$ws.Range("H107").Value = 3693.7778
$ws.Range("I107").Value = 1617.6
$ws.Range("J107").Value = 6289
$ws.Range("K107").Value = 1617.6
$ws.Range("L107").Value = 6289
$ws.Range("M107").Value = 302.4000000000001
$ws.Range("N107").Value = -10129
$ws.Range("H128").Value = 80779
$ws.Range("J128").Value = 80779
$ws.Range("L128").Value = 80779
$ws.Range("N128").Value = -90739
$ws.Range("H132").Value = 3763
$ws.Range("I132").Value = 3626.5
$ws.Range("J132").Value = 3899.5
$ws.Range("K132").Value = 10879.5
$ws.Range("L132").Value = 11698.5
$ws.Range("M132").Value = -8349.5
$ws.Range("N132").Value = -16758.5
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 975.5
$ws.Range("I16").Value = 975.5
$ws.Range("K16").Value = 975.5
$ws.Range("M16").Value = -805.5
$ws.Range("H82").Value = 34831.668
$ws.Range("I82").Value = 23496.75
$ws.Range("K82").Value = 23496.75
$ws.Range("M82").Value = -23135.75
$ws.Range("H85").Value = 34831.668
$ws.Range("I85").Value = 23496.75
$ws.Range("K85").Value = 23496.75
$ws.Range("M85").Value = -22248.75
$ws.Range("H136").Value = 2997.75
$ws.Range("I136").Value = 2997.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8993.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6443.25
$ws.Range("N136").ClearContents()
